$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 3, shifting existing rows 3-6 down to 4-7.
$ws.Rows.Item(3).Insert()
# The inserted row picks up formatting copied down from row 2; clear it so the
# row is truly blank (no leftover styled-but-empty cells), matching a plain
# "Insert Sheet Rows" on an otherwise-empty row.
$ws.Rows.Item(3).Clear()

# Row-insert doesn't auto-shift the hyperlink anchors in this engine, so
# rebuild the hyperlinks collection against the new (shifted) cell addresses.
# Each EMAIL_ADDRESS cell's mailto: hyperlink target matches its own text, so
# the correct targets can be recomputed straight from the (already-shifted)
# cell values instead of being hard-coded. Re-add in the same relative order
# as the original collection (the row that used to be C3 first, then the
# untouched C2, then the old C4/C5/C6 rows) so each link keeps its original
# relationship slot, just pointed at its new (shifted) cell.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C4"), ("mailto:" + $ws.Range("C4").Value()))
$ws.Hyperlinks.Add($ws.Range("C2"), ("mailto:" + $ws.Range("C2").Value()))
$ws.Hyperlinks.Add($ws.Range("C5"), ("mailto:" + $ws.Range("C5").Value()))
$ws.Hyperlinks.Add($ws.Range("C6"), ("mailto:" + $ws.Range("C6").Value()))
$ws.Hyperlinks.Add($ws.Range("C7"), ("mailto:" + $ws.Range("C7").Value()))

# Adding a hyperlink re-applies the "Hyperlink" cell style via a freshly
# created style index; restore the original shared style on each cell so the
# cell formatting matches what it was before the rebuild.
$ws.Range("C2").Style = "Hyperlink"
$ws.Range("C4").Style = "Hyperlink"
$ws.Range("C5").Style = "Hyperlink"
$ws.Range("C6").Style = "Hyperlink"
$ws.Range("C7").Style = "Hyperlink"

# Update the selection to match the post-edit state.
$ws.Range("C12").Select()
